$d = $word.ActiveDocument

# Locate the "© 2020 ... Powered by Jekyll ..." copyright paragraph that
# used to sit right before the final empty / page-break paragraphs, near
# the end of the document (right after the "Requisitos" section).
$targetIndex = -1
$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 1; $i--) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -like "*Powered by Jekyll*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -gt 0) {
    # Remove, from last to first so indices of earlier paragraphs remain
    # valid: the copyright paragraph itself, the page-break paragraph
    # right before it, and the empty paragraph before that.
    $d.Paragraphs($targetIndex).Range.Delete()
    $d.Paragraphs($targetIndex - 1).Range.Delete()
    $d.Paragraphs($targetIndex - 2).Range.Delete()
}
